# RW-3 Connector Master List.xlsx update
# - Remove the old "Mid MCN breakout" test connectors X-4000 / X-4001 from the
#   Connectors sheet (rows 71-72), since the Dyno System Diagram (v2.1) no
#   longer uses them.
# - Add the new Dyno MCN breakout connectors X-5002 / X-5003 to the
#   Connectors sheet.
# - Add the matching new cable assembly C-1006 to the Cable Assemblies sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Connectors"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Connectors")

# Remove the retired X-4000 / X-4001 rows (just clear their contents - the
# rows below do NOT shift up, row 73 was already a blank spacer row).
$ws1.Rows("71:72").ClearContents()

# Make room for two new connector rows right after the existing X-5001 row
# (row 93), pushing the "Anderson" section (and everything after it) down.
$ws1.Rows("94:96").Insert()

$ws1.Range("A94").Value2 = "X-5002"
$ws1.Range("B94").Value2 = "10-crkt MicroClasp"
$ws1.Range("C94").Value2 = "Dyno Electrical System Diagram"
$ws1.Range("E94").Value2 = "C-1006"

$ws1.Range("A95").Value2 = "X-5003"
$ws1.Range("B95").Value2 = "10-crkt MicroClasp"
$ws1.Range("C95").Value2 = "Dyno Electrical System Diagram"
$ws1.Range("E95").Value2 = "C-1000"

$ws1.Range("D90").Select()

# ---------------------------------------------------------------------
# Sheet "Cable Assemblies"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Cable Assemblies")

$ws2.Range("A45").Value2 = "C-1006"
$ws2.Range("B45").Value2 = 9
$ws2.Range("F45").Value2 = "Dyno Flow Sensor & Temp Sensors to Dyno MCN Backplane"
$ws2.Rows(45).RowHeight = 45

$ws2.Range("I45").Select()

# Leave "Connectors" as the active sheet/tab, matching the original workbook.
$ws1.Activate()
